$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.157.54'
$ws.Range("E2").Value = '  -0.27%  '

$ws.Range("D3").Value = '3.886.87'
$ws.Range("E3").Value = '  -0.33%  '

$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '482.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.79%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.623'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.84%  '

$ws.Range("E8").Value = '  -0.09%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.743'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.46%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.181'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000355'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.94%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.08'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.93%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.51'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.86%  '

$ws.Range("D14").Value = '4.506.35'
$ws.Range("E14").Value = '  -0.50%  '

$ws.Range("D15").Value = '3.881.38'
$ws.Range("E15").Value = '  -0.68%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.25'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.77%  '

$ws.Range("E17").Value = '  -0.50%  '

$ws.Range("E18").Value = '  +1.66%  '

$ws.Range("E19").Value = '  +0.93%  '

$ws.Range("D20").Value = '68.237.24'
$ws.Range("E20").Value = '  -0.38%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '429.26'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.22%  '

$ws.Range("E22").Value = '  +8.34%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.79'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.23%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.67'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +20.09%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '88.80'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.56%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.67'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.78%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.99%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.23'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.53%  '

$ws.Range("E29").Value = '  -2.85%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '720.30'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.49'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.42%  '

$ws.Range("E32").Value = '  +0.61%  '

$ws.Range("E33").Value = '  +3.27%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '61.78'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.56%  '

$ws.Range("B35").Value = 'PEPE'
$ws.Range("C35").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D35").Value = '0.0₃0879'
$ws.Range("E35").Value = '  -0.99%  '

$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.06'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.59%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '40.81'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.91%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.400'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +17.80%  '

$ws.Range("E39").Value = '  -2.87%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.00%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0497'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.26%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.70%  '

$ws.Range("E43").Value = '  +3.68%  '

$ws.Range("E44").Value = '  -3.21%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.38'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.47%  '

$ws.Range("E47").Value = '  -0.03%  '

$ws.Range("D48").Value = '0.0₆0359'
$ws.Range("E48").Value = '  +31.44%  '

$ws.Range("E49").Value = '  -0.66%  '

$ws.Range("E50").Value = '  -1.70%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '144.49'
$ws.Range("D51").Style = "Normal"
